$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "dnsfqnfkjdn"
$ws.Range("B1").Value = "flkfndslkfn"
$ws.Range("C1").Value = " "

$ws.Range("C1").Select()
